# Auto-generated edits reproducing the commit diff for Hyperion_Profits workbook
# (FFXIV Hyperion-server leve crafting-profit tracker: columns H..N are market-price /
# leve-profit figures recomputed from updated Universalis price data.)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: "One for the Road" / "Potion"
$ws.Range("H17").Value = 5320.2363
$ws.Range("J17").Value = 5320.2363
$ws.Range("L17").Value = 15960.7089
$ws.Range("N17").Value = -16296.7089

# Row 39: "Riches' Brew" / "Hi-Potion of Mind"
$ws.Range("H39").Value = 89.82353000000001
$ws.Range("I39").Value = 73.181816
$ws.Range("J39").Value = 120.333336
$ws.Range("K39").Value = 219.545448
$ws.Range("L39").Value = 361.000008
$ws.Range("M39").Value = 76.45455200000001
$ws.Range("N39").Value = -953.000008

# Row 53: "No Accounting for Waste" / "Enchanted Electrum Ink"
$ws.Range("H53").Value = 12752.25
$ws.Range("I53").Value = 264.3
$ws.Range("K53").Value = 264.3
$ws.Range("M53").Value = 372.7

# Row 74: "Adhesive of Antipathy" / "Wing Glue"
$ws.Range("H74").Value = 7221.4614
$ws.Range("I74").Value = 3633
$ws.Range("J74").Value = 7689.522
$ws.Range("K74").Value = 3633
$ws.Range("L74").Value = 7689.522
$ws.Range("M74").Value = -2697
$ws.Range("N74").Value = -9561.522000000001

# Row 77: "It's Gonna Grow Back (L)" / "Wing Glue"
$ws.Range("H77").Value = 7221.4614
$ws.Range("I77").Value = 3633
$ws.Range("J77").Value = 7689.522
$ws.Range("K77").Value = 18165
$ws.Range("L77").Value = 38447.61
$ws.Range("M77").Value = -13485
$ws.Range("N77").Value = -47807.61

# Row 98: "The Dotted Line" / "Enchanted Durium Ink"
$ws.Range("H98").Value = 1197.12
$ws.Range("I98").Value = 997
$ws.Range("K98").Value = 997
$ws.Range("M98").Value = 501

# Row 99: "Rumor Has It" / "Commanding Craftsman's Tea"
$ws.Range("H99").Value = 299
$ws.Range("I99").Value = 242.1
$ws.Range("J99").Value = 868
$ws.Range("K99").Value = 726.3
$ws.Range("L99").Value = 2604
$ws.Range("M99").Value = 771.7
$ws.Range("N99").Value = -5600

# Row 100: "Asking for a Friend" / "Beetle Glue"
$ws.Range("H100").Value = 2839
$ws.Range("I100").Value = 2686.8462
$ws.Range("K100").Value = 2686.8462
$ws.Range("M100").Value = -2145.8462

# Row 122: "Wishful Inking" / "Enchanted High Durium Ink"
$ws.Range("H122").Value = 1197.12
$ws.Range("I122").Value = 997
$ws.Range("K122").Value = 2991
$ws.Range("M122").Value = -541

# Row 137: "Cutting Edge of Culinary Quality" / "Magnesia Whetstone"
$ws.Range("H137").Value = 2976.8948
$ws.Range("I137").Value = 2445.3809
$ws.Range("K137").Value = 7336.1427
$ws.Range("M137").Value = -4786.1427

# Row 138: "All-night Crafting" / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 3036.0505
$ws.Range("I138").Value = 1713.625
$ws.Range("J138").Value = 3667.6567
$ws.Range("K138").Value = 5140.875
$ws.Range("L138").Value = 11002.9701
$ws.Range("M138").Value = -0.875
$ws.Range("N138").Value = -21282.9701

$ws = $wb.Worksheets.Item("ARM")
# Row 45: "Hollow Hallmarks" / "Mythril Ingot"
$ws.Range("H45").Value = 6853904.5
$ws.Range("J45").Value = 10593.4
$ws.Range("L45").Value = 10593.4
$ws.Range("N45").Value = -11347.4

# Row 61: "Dealing with the Tough Stuff" / "Cobalt Ingot"
$ws.Range("H61").Value = 4096.4634
$ws.Range("I61").Value = 4188.421
$ws.Range("K61").Value = 4188.421
$ws.Range("M61").Value = -3976.421

# Row 74: "As the Bolt Flies" / "Titanium Nugget"
$ws.Range("H74").Value = 312813.44
$ws.Range("I74").Value = 216801.33
$ws.Range("K74").Value = 216801.33
$ws.Range("M74").Value = -215927.33

# Row 77: "Heavy Metal Banned (L)" / "Titanium Nugget"
$ws.Range("H77").Value = 312813.44
$ws.Range("I77").Value = 216801.33
$ws.Range("K77").Value = 1084006.65
$ws.Range("M77").Value = -1079638.65

# Row 122: "Haste for High Durium" / "High Durium Nugget"
$ws.Range("H122").Value = 2607785
$ws.Range("I122").Value = 4281.4
$ws.Range("J122").Value = 6946957.5
$ws.Range("K122").Value = 12844.2
$ws.Range("L122").Value = 20840872.5
$ws.Range("M122").Value = -10394.2
$ws.Range("N122").Value = -20845772.5

# Row 127: "Once and for Alchemy" / "Bismuth Alembic"
$ws.Range("H127").Value = 32890
$ws.Range("J127").Value = 32890
$ws.Range("L127").Value = 32890
$ws.Range("N127").Value = -42810

# Row 136: "Metal with Mettle" / "Cobalt Tungsten Ingot"
$ws.Range("H136").Value = 4096.4634
$ws.Range("I136").Value = 4188.421
$ws.Range("K136").Value = 12565.263
$ws.Range("M136").Value = -10015.263

$ws = $wb.Worksheets.Item("BSM")
# Row 96: "Hammer Time" / "High Steel Sledgehammer"
$ws.Range("H96").Value = 17184
$ws.Range("I96").Value = 10047.833
$ws.Range("J96").Value = 60001
$ws.Range("K96").Value = 10047.833
$ws.Range("L96").Value = 60001
$ws.Range("M96").Value = -7301.833000000001
$ws.Range("N96").Value = -65493

# Row 104: "Hammer and Sails" / "Molybdenum Ball-pein Hammer"
$ws.Range("H104").Value = 30240
$ws.Range("J104").Value = 30240
$ws.Range("L104").Value = 30240
$ws.Range("N104").Value = -37228

$ws = $wb.Worksheets.Item("CRP")
# Row 31: "Wall Not Found" / "Walnut Lumber"
$ws.Range("H31").Value = 6605.6665
$ws.Range("I31").Value = 11438
$ws.Range("K31").Value = 11438
$ws.Range("M31").Value = -11143

# Row 34: "Armoires of the Rich and Famous" / "Walnut Lumber"
$ws.Range("H34").Value = 6605.6665
$ws.Range("I34").Value = 11438
$ws.Range("K34").Value = 11438
$ws.Range("M34").Value = -11236

# Row 58: "You Do the Heavy Lifting" / "Mahogany Lumber"
$ws.Range("H58").Value = 1830.1904
$ws.Range("I58").Value = 1521.0714
$ws.Range("K58").Value = 1521.0714
$ws.Range("M58").Value = -1318.0714

# Row 94: "Beech, Please" / "Beech Lumber"
$ws.Range("H94").Value = 1373.4445
$ws.Range("I94").Value = 806
$ws.Range("K94").Value = 806
$ws.Range("M94").Value = -355

# Row 132: "Hull Lotta Damage" / "Ginseng Lumber"
$ws.Range("H132").Value = 64981.375
$ws.Range("J132").Value = 1498.5
$ws.Range("L132").Value = 4495.5
$ws.Range("N132").Value = -9555.5

# Row 134: "Wood You Be Quiet" / "Ceiba Lumber"
$ws.Range("H134").Value = 25978.104
$ws.Range("I134").Value = 33677.562
$ws.Range("K134").Value = 101032.686
$ws.Range("M134").Value = -98497.68599999999

# Row 136: "Turali Quality" / "Dark Mahogany Lumber"
$ws.Range("H136").Value = 1830.1904
$ws.Range("I136").Value = 1521.0714
$ws.Range("K136").Value = 4563.2142
$ws.Range("M136").Value = -2013.2142

$ws = $wb.Worksheets.Item("CUL")
# Row 8: "Whip It" / "Sweet Cream"
$ws.Range("H8").Value = 255.09091
$ws.Range("I8").Value = 255.09091
$ws.Range("K8").Value = 765.27273
$ws.Range("M8").Value = -626.27273

# Row 38: "Pretty as a Picture" / "Dark Vinegar"
$ws.Range("H38").Value = 113
$ws.Range("J38").Value = 223.5
$ws.Range("L38").Value = 670.5
$ws.Range("N38").Value = -1364.5

# Row 56: "Culture Club" / "Crowned Pie"
$ws.Range("H56").Value = 17863424
$ws.Range("I56").Value = 17863424
$ws.Range("K56").Value = 17863424
$ws.Range("M56").Value = -17862894

# Row 122: "Salt of the North" / "Northern Sea Salt"
$ws.Range("H122").Value = 713.4167
$ws.Range("J122").Value = 646.5
$ws.Range("L122").Value = 5818.5
$ws.Range("N122").Value = -10718.5

# Row 129: "Comfort Food" / "Yakow Moussaka"
$ws.Range("H129").Value = 1084.9286
$ws.Range("I129").Value = 836.36365
$ws.Range("J129").Value = 1996.3334
$ws.Range("K129").Value = 2509.09095
$ws.Range("L129").Value = 5989.0002
$ws.Range("M129").Value = 2490.90905
$ws.Range("N129").Value = -15989.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 80: "Needs More Prayerbell" / "Hardsilver Ingot"
$ws.Range("H80").Value = 1886171.2
$ws.Range("I80").Value = 4083324
$ws.Range("K80").Value = 4083324
$ws.Range("M80").Value = -4082326

# Row 83: "With a Noise That Reaches Heaven (L)" / "Hardsilver Ingot"
$ws.Range("H83").Value = 1886171.2
$ws.Range("I83").Value = 4083324
$ws.Range("K83").Value = 20416620
$ws.Range("M83").Value = -20411628

# Row 102: "Put the Metal to the Peddle" / "Durium Ingot"
$ws.Range("H102").Value = 9064525
$ws.Range("I102").Value = 15875339
$ws.Range("K102").Value = 15875339
$ws.Range("M102").Value = -15873717

# Row 103: "Ring in the New" / "Azurite Ring of Fending"
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

# Row 113: "Copious Crystal Cannons" / "Manasilver Nugget"
$ws.Range("H113").Value = 27780210
$ws.Range("I113").Value = 55557420
$ws.Range("K113").Value = 55557420
$ws.Range("M113").Value = -55555250

# Row 138: "Orders Anonymous" / "White Gold Halfmask of Maiming"
$ws.Range("H138").Value = 67499.5
$ws.Range("J138").Value = 64999
$ws.Range("L138").Value = 64999
$ws.Range("N138").Value = -75279

$ws = $wb.Worksheets.Item("LTW")
# Row 61: "Spelling Me Softly" / "Raptor Leather"
$ws.Range("H61").Value = 11116917
$ws.Range("I61").Value = 13894750
$ws.Range("J61").Value = 5584.5
$ws.Range("K61").Value = 13894750
$ws.Range("L61").Value = 5584.5
$ws.Range("M61").Value = -13894548
$ws.Range("N61").Value = -5988.5

# Row 93: "Hide to Go Seek" / "Gagana Leather"
$ws.Range("H93").Value = 9261723
$ws.Range("I93").Value = 12347918
$ws.Range("J93").Value = 3140.5557
$ws.Range("K93").Value = 12347918
$ws.Range("L93").Value = 3140.5557
$ws.Range("M93").Value = -12346670
$ws.Range("N93").Value = -5636.5557

# Row 113: "Peace in Rest" / "Atrociraptor Leather"
$ws.Range("H113").Value = 11116917
$ws.Range("I113").Value = 13894750
$ws.Range("J113").Value = 5584.5
$ws.Range("K113").Value = 13894750
$ws.Range("L113").Value = 5584.5
$ws.Range("M113").Value = -13892580
$ws.Range("N113").Value = -9924.5

# Row 122: "Hell on Leather" / "Gaja Leather"
$ws.Range("H122").Value = 8221
$ws.Range("J122").Value = 9296.333000000001
$ws.Range("L122").Value = 27888.999
$ws.Range("N122").Value = -32788.999

# Row 132: "Tenets of Tanning" / "Silver Lobo Leather"
$ws.Range("H132").Value = 8515.438
$ws.Range("I132").Value = 8655.854499999999
$ws.Range("K132").Value = 25967.5635
$ws.Range("M132").Value = -23437.5635

$ws = $wb.Worksheets.Item("WVR")
# Row 41: "Half Is the New Double" / "Linen Halfgloves"
$ws.Range("H41").Value = 89823.336
$ws.Range("J41").Value = 89823.336
$ws.Range("L41").Value = 89823.336
$ws.Range("N41").Value = -90603.336

# Row 132: "Comfy Cabins" / "Snow Cotton Cloth"
$ws.Range("H132").Value = 25261290
$ws.Range("I132").Value = 27787184
$ws.Range("K132").Value = 83361552
$ws.Range("M132").Value = -83359022

